$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "62.582.56"
Set-TextValue $ws.Range("E2") "  +2.60%  "
Set-TextValue $ws.Range("D3") "2.949.75"
Set-TextValue $ws.Range("E3") "  +2.38%  "
Set-TextValue $ws.Range("E4") "  -0.12%  "
Set-TextValue $ws.Range("D5") "588.41"
Set-TextValue $ws.Range("E5") "  +0.14%  "
Set-TextValue $ws.Range("D6") "146.67"
Set-TextValue $ws.Range("E6") "  +4.62%  "
Set-TextValue $ws.Range("E7") "  -0.11%  "
Set-TextValue $ws.Range("D8") "2.949.21"
Set-TextValue $ws.Range("E8") "  +2.53%  "
Set-TextValue $ws.Range("E9") "  +2.96%  "
Set-TextValue $ws.Range("D10") "6.97"
Set-TextValue $ws.Range("E10") "  +1.51%  "
Set-TextValue $ws.Range("E11") "  +9.06%  "
Set-TextValue $ws.Range("E12") "  +1.72%  "
Set-TextValue $ws.Range("E13") "  +6.97%  "
Set-TextValue $ws.Range("D14") "32.17"
Set-TextValue $ws.Range("E14") "  -0.29%  "
Set-TextValue $ws.Range("E15") "  -1.03%  "
Set-TextValue $ws.Range("D16") "3.436.71"
Set-TextValue $ws.Range("E16") "  +2.12%  "
Set-TextValue $ws.Range("D17") "62.551.79"
Set-TextValue $ws.Range("E17") "  +2.46%  "
Set-TextValue $ws.Range("D18") "2.956.20"
Set-TextValue $ws.Range("E18") "  +2.20%  "
Set-TextValue $ws.Range("E19") "  +2.35%  "
Set-TextValue $ws.Range("D20") "434.01"
Set-TextValue $ws.Range("E20") "  +1.90%  "
Set-TextValue $ws.Range("D21") "13.46"
Set-TextValue $ws.Range("E21") "  +1.62%  "
Set-TextValue $ws.Range("D22") "0.659"
Set-TextValue $ws.Range("E22") "  +1.28%  "
Set-TextValue $ws.Range("E23") "  +0.64%  "
Set-TextValue $ws.Range("E24") "  +6.34%  "
Set-TextValue $ws.Range("D25") "80.12"
Set-TextValue $ws.Range("E25") "  +0.35%  "
Set-TextValue $ws.Range("D26") "11.87"
Set-TextValue $ws.Range("E26") "  +4.80%  "
Set-TextValue $ws.Range("E27") "  +2.31%  "
Set-TextValue $ws.Range("E28") "  -0.03%  "
Set-TextValue $ws.Range("D29") "7.16"
Set-TextValue $ws.Range("E29") "  +6.89%  "
Set-TextValue $ws.Range("D30") "2.17"
Set-TextValue $ws.Range("E30") "  +3.82%  "
Set-TextValue $ws.Range("D31") "2.57"
Set-TextValue $ws.Range("E31") "  +1.84%  "
Set-TextValue $ws.Range("E32") "  +17.89%  "
Set-TextValue $ws.Range("E33") "  +3.16%  "
Set-TextValue $ws.Range("D34") "26.16"
Set-TextValue $ws.Range("E34") "  +1.47%  "
Set-TextValue $ws.Range("E35") "  -0.10%  "
Set-TextValue $ws.Range("D36") "0.989"
Set-TextValue $ws.Range("E36") "  +2.02%  "
Set-TextValue $ws.Range("E37") "  +2.34%  "
Set-TextValue $ws.Range("E38") "  +7.06%  "
Set-TextValue $ws.Range("D39") "49.62"
Set-TextValue $ws.Range("E39") "  +1.12%  "
Set-TextValue $ws.Range("E40") "  +4.71%  "
Set-TextValue $ws.Range("E41") "  +0.41%  "
Set-TextValue $ws.Range("D42") "0.115"
Set-TextValue $ws.Range("E42") "  -2.00%  "
Set-TextValue $ws.Range("E43") "  +3.72%  "
Set-TextValue $ws.Range("D44") "39.18"
Set-TextValue $ws.Range("E44") "  -2.82%  "
Set-TextValue $ws.Range("D45") "134.94"
Set-TextValue $ws.Range("E45") "  +1.70%  "
Set-TextValue $ws.Range("D46") "2.676.83"
Set-TextValue $ws.Range("E46") "  +0.89%  "
Set-TextValue $ws.Range("E47") "  +0.27%  "
Set-TextValue $ws.Range("D48") "353.18"
Set-TextValue $ws.Range("E48") "  +2.98%  "
Set-TextValue $ws.Range("E50") "  +1.47%  "
Set-TextValue $ws.Range("D51") "22.49"
Set-TextValue $ws.Range("E51") "  -0.19%  "
